$d = $word.ActiveDocument

# Locate the paragraph that ends the "29.09.2023" entry (the one whose last
# run is a manual line break <w:br/>, right after the sentence about
# realizing the sass job mistake). In the original document it is directly
# followed by an already-existing empty BodyText paragraph and then the
# page-break paragraph / sectPr.
$anchorIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.Contains("15-20 minutes until I realized")) {
        $anchorIdx = $i
        break
    }
}

if ($anchorIdx -eq -1) {
    throw "Could not locate anchor paragraph for the 30.09.2023 diary entry"
}

$anchor = $d.Paragraphs.Item($anchorIdx)

# Create one new empty paragraph right after it; InsertXML below expands it
# into the three paragraphs that belong to the new "30.09.2023" diary entry.
$anchor.Range.InsertParagraphAfter()
$target = $d.Paragraphs.Item($anchorIdx + 1)
$r = $target.Range

$apos = [char]0x2019

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p><w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>30</w:t></w:r>' +
'<w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>.09.2023</w:t></w:r>' +
'</w:p>' +
'<w:p><w:pPr><w:pStyle w:val="BodyText"/><w:jc w:val="left"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
('<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>I finished the third video where I' + $apos + 've learned how to</w:t></w:r>') +
'<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> implement changes onto the html class dynamically using </w:t></w:r>' +
'<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>javascript</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
'<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">. In this case the idea behind it all is to change the class </w:t></w:r>' +
'<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">of the element on click so the </w:t></w:r>' +
'<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>css</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
'<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> format changes and adapts to </w:t></w:r>' +
'<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>that different class.</w:t></w:r>' +
'</w:p>' +
'<w:p><w:pPr><w:pStyle w:val="BodyText"/><w:jc w:val="left"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Also, I learned how to use the transformation </w:t></w:r>' +
'<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>different functions to relocate an element on the website (the menu button)</w:t></w:r>' +
'</w:p>' +
'</w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

[void]$r.InsertXML($xml)

Write-Output "Inserted the 30.09.2023 diary entry after paragraph $anchorIdx"
